# Fix sum response from import dssv:
# The "Email" column (J) is removed from the error-report sheet. Columns
# that were to its right (Quê quán, Ghi chú) shift left by one. Row 13 is
# special-cased: its old J13 already held the real error text ("Quê quán
# không hợp lệ") while K13 held stray leftover data ("Ái chà!!!") that
# should not survive the shift - after the fix J13 keeps its original
# error text/highlight and K13 ends up blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 is the one row where the naive "shift everything left" rule
# does not apply to column J - capture its original value/highlight so
# we can restore them after the column delete.
$keepJ13Value = $ws.Range("J13").Value2
$keepJ13ColorIndex = $ws.Range("J13").Interior.ColorIndex

# Deleting the entire "Email" column shifts Quê quán/Ghi chú (and their
# cell formatting) one column to the left for every row.
$ws.Columns("J").Delete()

# Undo the unwanted shift for row 13 only. (K13 already ends up blank on
# its own, because the old L13 - which shifts into K13 - was blank too.)
$ws.Range("J13").Value = $keepJ13Value
$ws.Range("J13").Interior.ColorIndex = $keepJ13ColorIndex
